$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "test cases"
$ws2 = $wb.Worksheets.Item(2)   # "settings"

# Rename the fee-type related shared strings used on the "test cases" sheet:
#   "Kazakhstan Country" -> "Kazakhstan FEE"
#   "KZC"                -> "KZF"
#   "USA Country"        -> "USA FEE"
# Every cell that referenced the old strings is rewritten so the old shared
# strings become unused and get dropped when the workbook is saved, while the
# new strings are appended to the shared string table.
$ws1.Range("A2").Value = "Kazakhstan FEE"
$ws1.Range("B2").Value = "KZF"
$ws1.Range("B3").Value = "KZF"
$ws1.Range("A4").Value = "Kazakhstan FEE"
$ws1.Range("B4").Value = "KZF"
$ws1.Range("A6").Value = "USA FEE"

# Update the "settings" sheet's selection first, while it is still the
# active sheet (selecting a range on the other sheet re-activates it, so this
# must happen before we switch sheets below).
$ws2.Range("A1:B4").Select()

# Switch the active tab from "settings" to "test cases" and move the
# selection there, matching the new view state captured in the workbook.
$ws1.Activate()
$ws1.Range("A6").Select()
